# Add a new "Croatia" Test Data sheet, based on the existing "Turkey" sheet,
# matching the pattern used for all the other per-country sheets in this
# workbook (Germany, Belgium, Czech, Swiss, Portugal, Slovakia, Italy, Spain,
# Turkey).

$wb = $excel.ActiveWorkbook

# Select the whole "Turkey" sheet (this is what makes Excel record the
# source sheet's new "select all" selection state once we move away from it)
# and copy it to a new tab placed right after it.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Range("A1:XFD1048576").Select()
$turkey.Copy([System.Reflection.Missing]::Value, $turkey)

# The copy gets inserted immediately after "Turkey".
$newSheet = $wb.Worksheets.Item($turkey.Index + 1)
$newSheet.Name = "Croatia"

# Fill in the market-specific values for the new Croatia sheet.
$newSheet.Range("B4").Value = "NGC-3139/T2473"
$newSheet.Range("B2").Value = "Croatia Market"

# Leave the selection on the new sheet where it ended up.
$newSheet.Range("I13").Select()
